$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old row 3 (with text "333"/"321" in A3/B3 and 8/8 in C3/D3) is removed;
# row 2 now carries new text values ("456"/"123") plus new numeric results (10/10).
$ws.Rows.Item(3).Delete()

# Force A2/B2 to be stored as text (shared strings), not numbers, then drop the
# temporary number-format override so no extra style sticks to the cells.
$a2 = $ws.Cells.Item(2, 1)
$a2.NumberFormat = "@"
$a2.Value = "456"
$a2.Style = "Normal"

$b2 = $ws.Cells.Item(2, 2)
$b2.NumberFormat = "@"
$b2.Value = "123"
$b2.Style = "Normal"

$ws.Cells.Item(2, 3).Value = 10
$ws.Cells.Item(2, 4).Value = 10
